$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2291666666666667
$ws.Range("C2").Value = 0.4791666666666667
$ws.Range("J2").Value = 0.01488095238095238
$ws.Range("P2").Value = 0.1607142857142857
$ws.Range("S2").Value = 0.1160714285714286
$ws.Range("B3").Value = 0.005681818181818182
$ws.Range("C3").Value = 0.06818181818181818
$ws.Range("J3").Value = 0.02272727272727273
$ws.Range("P3").Value = 0.7159090909090909
$ws.Range("S3").Value = 0.1875
$ws.Range("J4").Value = 0.03448275862068965
$ws.Range("P4").Value = 0.5344827586206896
$ws.Range("S4").Value = 0.4310344827586207
$ws.Range("B6").Value = 0.0321285140562249
$ws.Range("D6").Value = 0.01606425702811245
$ws.Range("F6").Value = 0.06827309236947791
$ws.Range("J6").Value = 0.3132530120481928
$ws.Range("O6").Value = 0.02409638554216868
$ws.Range("Q6").Value = 0.1686746987951807
$ws.Range("R6").Value = 0.04819277108433735
$ws.Range("S6").Value = 0.3293172690763052
$ws.Range("B7").Value = 0.1012658227848101
$ws.Range("D7").Value = 0.01265822784810127
$ws.Range("F7").Value = 0.0379746835443038
$ws.Range("J7").Value = 0.1476793248945148
$ws.Range("O7").Value = 0.02953586497890295
$ws.Range("Q7").Value = 0.1687763713080169
$ws.Range("R7").Value = 0.0759493670886076
$ws.Range("S7").Value = 0.4261603375527426
$ws.Range("B8").Value = 0.08990825688073395
$ws.Range("D8").Value = 0.01834862385321101
$ws.Range("E8").Value = 0.001834862385321101
$ws.Range("F8").Value = 0.05137614678899083
$ws.Range("J8").Value = 0.06972477064220184
$ws.Range("O8").Value = 0.01467889908256881
$ws.Range("Q8").Value = 0.1834862385321101
$ws.Range("R8").Value = 0.108256880733945
$ws.Range("S8").Value = 0.4623853211009175
$ws.Range("B9").Value = 0.1040723981900453
$ws.Range("D9").Value = 0.01809954751131222
$ws.Range("F9").Value = 0.08144796380090498
$ws.Range("J9").Value = 0.05882352941176471
$ws.Range("O9").Value = 0.004524886877828055
$ws.Range("Q9").Value = 0.2262443438914027
$ws.Range("R9").Value = 0.06334841628959276
$ws.Range("S9").Value = 0.4434389140271493
$ws.Range("B10").Value = 0.09928151534944481
$ws.Range("D10").Value = 0.0248203788373612
$ws.Range("E10").Value = 0.0006531678641410843
$ws.Range("F10").Value = 0.05290659699542782
$ws.Range("J10").Value = 0.08556499020248204
$ws.Range("O10").Value = 0.0124101894186806
$ws.Range("Q10").Value = 0.2135858915741345
$ws.Range("R10").Value = 0.09601567602873938
$ws.Range("S10").Value = 0.4147615937295885
$ws.Range("G11").Value = 0.1770573566084788
$ws.Range("J11").Value = 0.07730673316708229
$ws.Range("K11").Value = 0.2394014962593516
$ws.Range("L11").Value = 0.4887780548628429
$ws.Range("S11").Value = 0.01745635910224439
$ws.Range("G12").Value = 0.6926605504587156
$ws.Range("J12").Value = 0.1422018348623853
$ws.Range("K12").Value = 0.02293577981651376
$ws.Range("L12").Value = 0.09174311926605505
$ws.Range("S12").Value = 0.05045871559633028
$ws.Range("G13").Value = 0.5777777777777777
$ws.Range("J13").Value = 0.3777777777777778
$ws.Range("S13").Value = 0.04444444444444445
$ws.Range("F15").Value = 0.04128440366972477
$ws.Range("H15").Value = 0.1376146788990826
$ws.Range("I15").Value = 0.07339449541284404
$ws.Range("J15").Value = 0.3853211009174312
$ws.Range("K15").Value = 0.08256880733944955
$ws.Range("M15").Value = 0.01376146788990826
$ws.Range("O15").Value = 0.03211009174311927
$ws.Range("S15").Value = 0.2339449541284404
$ws.Range("F16").Value = 0.02427184466019417
$ws.Range("H16").Value = 0.1504854368932039
$ws.Range("I16").Value = 0.06796116504854369
$ws.Range("J16").Value = 0.4368932038834951
$ws.Range("K16").Value = 0.1019417475728155
$ws.Range("M16").Value = 0.02427184466019417
$ws.Range("O16").Value = 0.05339805825242718
$ws.Range("S16").Value = 0.1407766990291262
$ws.Range("F17").Value = 0.02888086642599278
$ws.Range("H17").Value = 0.1823104693140794
$ws.Range("I17").Value = 0.0776173285198556
$ws.Range("J17").Value = 0.4512635379061372
$ws.Range("K17").Value = 0.06678700361010831
$ws.Range("M17").Value = 0.02527075812274368
$ws.Range("N17").Value = 0.001805054151624549
$ws.Range("O17").Value = 0.03971119133574007
$ws.Range("S17").Value = 0.1263537906137184
$ws.Range("F18").Value = 0.02016129032258064
$ws.Range("H18").Value = 0.1935483870967742
$ws.Range("I18").Value = 0.07661290322580645
$ws.Range("J18").Value = 0.4556451612903226
$ws.Range("K18").Value = 0.0846774193548387
$ws.Range("M18").Value = 0.02016129032258064
$ws.Range("O18").Value = 0.04838709677419355
$ws.Range("S18").Value = 0.1008064516129032
$ws.Range("F19").Value = 0.02088305489260143
$ws.Range("H19").Value = 0.1992840095465394
$ws.Range("J19").Value = 0.3723150357995227
$ws.Range("K19").Value = 0.1193317422434367
$ws.Range("M19").Value = 0.01312649164677804
$ws.Range("O19").Value = 0.05369928400954654
$ws.Range("S19").Value = 0.1426014319809069
